# Update "苏州-漫展信息" workbook: refresh scraped 漫展 (con) listing data.
# Applies to both the "展览" sheet and the "全部类型" sheet:
#   - 5 stale rows (old rows 3-7: 星河璀璨..., 第三届OCG Summer Carnival,
#     太仓第六届龙狮, 昆山第七届xcy新次元, OCG国潮凌飞内场) are removed -
#     they no longer appear in the refreshed feed.
#   - The remaining rows shift up and get refreshed "想去人数" (F column)
#     counts from the new scrape.
#   - Column A (sequence number) is renumbered 1..N to stay contiguous.

$wb = $excel.ActiveWorkbook

function Update-ConSheet {
    param($ws)

    # Remove the 5 rows that dropped out of the refreshed listing.
    $ws.Rows("3:7").Delete()

    # Renumber the sequence column (A) so it stays 1, 2, 3, ... with no gaps.
    $dims = $ws.UsedRange
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }
}

# --- Sheet "展览" ---------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
Update-ConSheet $wsExpo

$wsExpo.Range("F2").Value  = 278
$wsExpo.Range("F4").Value  = 1153
$wsExpo.Range("F5").Value  = 16530
$wsExpo.Range("F7").Value  = 1622
$wsExpo.Range("F10").Value = 200
$wsExpo.Range("F11").Value = 124
$wsExpo.Range("F12").Value = 11524
$wsExpo.Range("F13").Value = 23
$wsExpo.Range("F14").Value = 1186
$wsExpo.Range("F15").Value = 4556
$wsExpo.Range("F16").Value = 391
$wsExpo.Range("F17").Value = 398
$wsExpo.Range("F18").Value = 60
$wsExpo.Range("F19").Value = 867

# --- Sheet "全部类型" ------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
Update-ConSheet $wsAll

$wsAll.Range("F2").Value  = 278
$wsAll.Range("F5").Value  = 1153
$wsAll.Range("F6").Value  = 16530
$wsAll.Range("F8").Value  = 1622
$wsAll.Range("F11").Value = 200
$wsAll.Range("F12").Value = 124
$wsAll.Range("F15").Value = 11524
$wsAll.Range("F16").Value = 23
$wsAll.Range("F17").Value = 1186
$wsAll.Range("F18").Value = 4556
$wsAll.Range("F19").Value = 391
$wsAll.Range("F20").Value = 398
$wsAll.Range("F21").Value = 60
$wsAll.Range("F22").Value = 867
